$wb = $excel.ActiveWorkbook

# --- Rename the original sheet "Sheet 1" -> "datos" -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "datos"

# --- Add a new worksheet "metadatos" right after "datos" -------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadatos"

# --- Formatting: plain (non-scheme) Calibri 11 font on most of the table ---
# (cell A5 is intentionally left on the default style, matching source)
$ws2.Range("A1:D4").Font.Name = "Calibri"
$ws2.Range("A1:D4").Font.Size = 11
$ws2.Range("B5:D7").Font.Name = "Calibri"
$ws2.Range("B5:D7").Font.Size = 11
$ws2.Range("A6:A7").Font.Name = "Calibri"
$ws2.Range("A6:A7").Font.Size = 11

# --- Date number format on D2:D6, set before the values so the engine ----
# --- doesn't fabricate an intermediate auto-detected date format ---------
$ws2.Range("D2:D6").NumberFormat = "d-mmm-yy"

# --- Header row --------------------------------------------------------
$ws2.Range("A1").Value = "Variables"
$ws2.Range("B1").Value = "Descripción"
$ws2.Range("C1").Value = "Fuente"
$ws2.Range("D1").Value = "Fecha_de_extracción"

# --- Row 2: anno ---------------------------------------------------------
$ws2.Range("A2").Value = "anno"
$ws2.Range("B2").Value = "Año"
$ws2.Range("C2").Value = "…"

# --- Row 3: codmpio --------------------------------------------------------
$ws2.Range("A3").Value = "codmpio"
$ws2.Range("B3").Value = "Código del municipio"
$ws2.Range("C3").Value = "…"

# --- Row 4: SRPA_2 --------------------------------------------------------
$ws2.Range("A4").Value = "SRPA_2"
$ws2.Range("B4").Value = "No. de adolescentes que ingresan al SRPA con una medida no privativa de la libertad`nNo. de adolescentes que ingresan al SRPA con una medida privativa de la libertad "
$ws2.Range("C4").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- Row 5: ingresos_totales (A5 stays on the default/unstyled cell) -----
$ws2.Range("A5").Value = "ingresos_totales"
$ws2.Range("B5").Value = " No. total de adolescentes que han ingresado al sistema SRPA en el mismo periodo y territorio. x 100"
$ws2.Range("C5").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- Row 6: tasa -----------------------------------------------------------
$ws2.Range("A6").Value = "tasa"
$ws2.Range("C6").Value = "Elaboración Propia"

# --- D2:D6 extraction date: 2025-03-06, stored with no time component ----
$d = Get-Date -Year 2025 -Month 3 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws2.Range("D2:D6").Value = $d

# --- Keep rows 4 & 5 at the standard 15pt height even though the ---------
# --- multi-line text in column B would otherwise auto-expand them --------
$ws2.Rows("4:5").RowHeight = 15

# --- Column D width (Excel originally auto-fit this to the date strings) -
$ws2.Columns("D").ColumnWidth = 19

# --- View: "metadatos" becomes the active/front sheet, selection D2:D6 ---
$ws2.Activate() | Out-Null
$ws2.Range("D2:D6").Select() | Out-Null
